$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the date style (s="2") for the new row 53 by copying format from row 52
$ws.Range("A52").Copy() | Out-Null
$ws.Range("A53").PasteSpecial(-4122) | Out-Null

# Bulk-write the full data block (rows 2-53, columns A-E) using a 2D array assignment
$data = New-Object 'object[,]' 52,5
$data[0,0] = 39400
$data[0,1] = 2007
$data[0,2] = 1.144978573787081
$data[0,3] = 2008
$data[0,4] = 2.918510996763723
$data[1,0] = 39583
$data[1,1] = 2008
$data[1,2] = 3.09825757489699
$data[1,3] = 2009
$data[1,4] = 7.865470614547343
$data[2,0] = 39765
$data[2,1] = 2008
$data[2,2] = 1.381024225294869
$data[2,3] = 2009
$data[2,4] = 4.264380119800992
$data[3,0] = 39948
$data[3,1] = 2009
$data[3,2] = -1.791203563722299
$data[3,3] = 2010
$data[3,4] = -6.760862998203643
$data[4,0] = 40130
$data[4,1] = 2009
$data[4,2] = -0.8792832172735965
$data[4,3] = 2010
$data[4,4] = 3.109784435759599
$data[5,0] = 40310
$data[5,1] = 2010
$data[5,2] = 0.5799958470386946
$data[5,3] = 2011
$data[5,4] = 5.643342995751777
$data[6,0] = 40494
$data[6,1] = 2010
$data[6,2] = 0.9337833426867448
$data[6,3] = 2011
$data[6,4] = 2.730731696345146
$data[7,0] = 40676
$data[7,1] = 2011
$data[7,2] = 3.181454202131073
$data[7,3] = 2012
$data[7,4] = 4.38978860149748
$data[8,0] = 40862
$data[8,1] = 2011
$data[8,2] = 2.791140000794279
$data[8,3] = 2012
$data[8,4] = 1.683857142130885
$data[9,0] = 41044
$data[9,1] = 2012
$data[9,2] = 0.5930547804883668
$data[9,3] = 2013
$data[9,4] = -1.194610791900008
$data[10,0] = 41228
$data[10,1] = 2012
$data[10,2] = 0.4451370000809973
$data[10,3] = 2013
$data[10,4] = -0.2532347529486723
$data[11,0] = 41409
$data[11,1] = 2013
$data[11,2] = -0.3951783438669754
$data[11,3] = 2014
$data[11,4] = 0.03694906323863378
$data[12,0] = 41592
$data[12,1] = 2013
$data[12,2] = 0.2545814083968478
$data[12,3] = 2014
$data[12,4] = 2.553470871380537
$data[13,0] = 41774
$data[13,1] = 2014
$data[13,2] = 3.292216014290039
$data[13,3] = 2015
$data[13,4] = 7.617133650412211
$data[14,0] = 41957
$data[14,1] = 2014
$data[14,2] = 1.297015177357297
$data[14,3] = 2015
$data[14,4] = -0.434146007584113
$data[15,0] = 42137
$data[15,1] = 2015
$data[15,2] = 1.670328650030184
$data[15,3] = 2016
$data[15,4] = 2.037906845818616
$data[16,0] = 42321
$data[16,1] = 2015
$data[16,2] = 1.365576377841027
$data[16,3] = 2016
$data[16,4] = 2.383242923544526
$data[17,0] = 42503
$data[17,1] = 2016
$data[17,2] = 2.562791874943371
$data[17,3] = 2017
$data[17,4] = 3.265947405805814
$data[18,0] = 42689
$data[18,1] = 2016
$data[18,2] = 2.204449574611278
$data[18,3] = 2017
$data[18,4] = 1.688977015142101
$data[19,0] = 42867
$data[19,1] = 2017
$data[19,2] = 1.526411006965533
$data[19,3] = 2018
$data[19,4] = 0.6601843988560674
$data[20,0] = 43053
$data[20,1] = 2017
$data[20,2] = 2.18621550610123
$data[20,3] = 2018
$data[20,4] = 2.066615940231964
$data[21,0] = 43145
$data[21,1] = 2018
$data[21,2] = 1.297923389414657
$data[21,3] = 2019
$data[21,4] = 0.5447775838346658
$data[22,0] = 43235
$data[22,1] = 2018
$data[22,2] = 1.63465618619294
$data[22,3] = 2019
$data[22,4] = 1.551857746372698
$data[23,0] = 43326
$data[23,1] = 2018
$data[23,2] = 2.054458927584024
$data[23,3] = 2019
$data[23,4] = 3.441981941009353
$data[24,0] = 43418
$data[24,1] = 2018
$data[24,2] = 1.911050033324102
$data[24,3] = 2019
$data[24,4] = 3.0862758122153
$data[25,0] = 43510
$data[25,1] = 2019
$data[25,2] = 1.064009474888983
$data[25,3] = 2020
$data[25,4] = 0.03490120525229123
$data[26,0] = 43600
$data[26,1] = 2019
$data[26,2] = 1.35261353265177
$data[26,3] = 2020
$data[26,4] = 0.8024032016000104
$data[27,0] = 43691
$data[27,1] = 2019
$data[27,2] = 1.668617211002466
$data[27,3] = 2020
$data[27,4] = 1.816757311461803
$data[28,0] = 43783
$data[28,1] = 2019
$data[28,2] = 1.457852003181337
$data[28,3] = 2020
$data[28,4] = -1.135072001636328
$data[29,0] = 43875
$data[29,1] = 2020
$data[29,2] = 1.60064760240497
$data[29,3] = 2021
$data[29,4] = 2.904532120297287
$data[30,0] = 43966
$data[30,1] = 2020
$data[30,2] = -2.082763426755907
$data[30,3] = 2021
$data[30,4] = -5.866344937500023
$data[31,0] = 44068
$data[31,1] = 2020
$data[31,2] = -3.357986809108382
$data[31,3] = 2021
$data[31,4] = 11.62806235225531
$data[32,0] = 44159
$data[32,1] = 2020
$data[32,2] = -3.258619210312885
$data[32,3] = 2021
$data[32,4] = -2.878617960200258
$data[33,0] = 44251
$data[33,1] = 2021
$data[33,2] = -3.88825249955117
$data[33,3] = 2022
$data[33,4] = -7.923811177410267
$data[34,0] = 44341
$data[34,1] = 2021
$data[34,2] = -0.1380317107957718
$data[34,3] = 2022
$data[34,4] = 7.749494937649115
$data[35,0] = 44432
$data[35,1] = 2021
$data[35,2] = -0.1388955462784724
$data[35,3] = 2022
$data[35,4] = 4.532186626383039
$data[36,0] = 44525
$data[36,1] = 2021
$data[36,2] = 0.4255262881966981
$data[36,3] = 2022
$data[36,4] = 1.466936654457096
$data[37,0] = 44617
$data[37,1] = 2022
$data[37,2] = 3.299288015397095
$data[37,3] = 2023
$data[37,4] = 3.813885712818554
$data[38,0] = 44706
$data[38,1] = 2022
$data[38,2] = 3.848999231984762
$data[38,3] = 2023
$data[38,4] = 2.866003071127765
$data[39,0] = 44798
$data[39,1] = 2022
$data[39,2] = 3.463320568938566
$data[39,3] = 2023
$data[39,4] = 0.9167463358189964
$data[40,0] = 44890
$data[40,1] = 2022
$data[40,2] = 3.293290997728171
$data[40,3] = 2023
$data[40,4] = -1.421977974472588
$data[41,0] = 44981
$data[41,1] = 2023
$data[41,2] = -0.4618455958399603
$data[41,3] = 2024
$data[41,4] = -2.262006095280478
$data[42,0] = 45071
$data[42,1] = 2023
$data[42,2] = -0.3745803349312071
$data[42,3] = 2024
$data[42,4] = 0.645722451525943
$data[43,0] = 45163
$data[43,1] = 2023
$data[43,2] = -0.2871245688614854
$data[43,3] = 2024
$data[43,4] = -0.1008213866759977
$data[44,0] = 45254
$data[44,1] = 2023
$data[44,2] = -0.2814561130375925
$data[44,3] = 2024
$data[44,4] = -0.3873858053678236
$data[45,0] = 45345
$data[45,1] = 2024
$data[45,2] = -0.4631846496550684
$data[45,3] = 2025
$data[45,4] = 0.11654623407098
$data[46,0] = 45436
$data[46,1] = 2024
$data[46,2] = -1.091476630333243
$data[46,3] = 2025
$data[46,4] = -0.4907904687545206
$data[47,0] = 45534
$data[47,1] = 2024
$data[47,2] = -0.9636841177852018
$data[47,3] = 2025
$data[47,4] = -0.1102571493959759
$data[48,0] = 45618
$data[48,1] = 2024
$data[48,2] = -0.6470065423293869
$data[48,3] = 2025
$data[48,4] = 1.276847713071927
$data[49,0] = 45713
$data[49,1] = 2025
$data[49,2] = 2.749014873207956
$data[49,3] = 2026
$data[49,4] = 1.993341940459148
$data[50,0] = 45800
$data[50,1] = 2025
$data[50,2] = 1.626992717807862
$data[50,3] = 2026
$data[50,4] = 0.6270138473519316
$data[51,0] = 45891
$data[51,1] = 2025
$data[51,2] = 2.053865394798304
$data[51,3] = 2026
$data[51,4] = 0.8600861498751833

$ws.Range("A2:E53").Value = $data
